$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "backplane-breaout-big"

# Header row
$ws.Range("A1").Value = "Item"
$ws.Range("B1").Value = "Qty"
$ws.Range("C1").Value = "Reference(s)"
$ws.Range("D1").Value = "Value"
$ws.Range("E1").Value = "Footprint"
$ws.Range("F1").Value = "Manufacturer Part Number 1"

# Data rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 7
$ws.Range("C2").Value = "D1, D4, D5, D6, D9, D10, D13"
$ws.Range("D2").Value = "RED"
$ws.Range("E2").Value = "LED_SMD:LED_0603_1608Metric_Castellated"
$ws.Range("F2").Value = "LTST-C191KRKT"

$ws.Range("A3").Formula = "=A2+1"
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = "D2, D3, D7, D8, D11, D12, D14, D15"
$ws.Range("D3").Value = "GREEN"
$ws.Range("E3").Value = "LED_SMD:LED_0603_1608Metric_Castellated"
$ws.Range("F3").Value = "LTST-C191KGKT"

$ws.Range("A4").Formula = "=A3+1"
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = "D16, D19, D22, D23"
$ws.Range("D4").Value = "BLUE"
$ws.Range("E4").Value = "LED_SMD:LED_0603_1608Metric_Castellated"
$ws.Range("F4").Value = "LTST-C191TBKT"

$ws.Range("A5").Formula = "=A4+1"
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = "D17, D18, D20, D21"
$ws.Range("D5").Value = "ORANGE"
$ws.Range("E5").Value = "LED_SMD:LED_0603_1608Metric_Castellated"
$ws.Range("F5").Value = "LTST-C191KFKT"

$ws.Range("A6").Formula = "=A5+1"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = "H1, H2"
$ws.Range("D6").Value = "PC/104 Headers"
$ws.Range("E6").Value = "Connector_PinHeader_2.54mm:PinHeader_2x26_P2.54mm_Vertical"
$ws.Range("F6").Value = "ESQ-132-12-G-D"

$ws.Range("A7").Formula = "=A6+1"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = "J1, J2"
$ws.Range("D7").Value = "Breakout Connector"
$ws.Range("E7").Value = "backplane-breakout-big:SAMTEC_FTSH-130-04-L-DH_flipped"
$ws.Range("F7").Value = "FTSH-130-04-L-DH"

$ws.Range("A8").Formula = "=A7+1"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "J4, J5"
$ws.Range("D8").Value = "SSM-130-L-DV"
$ws.Range("E8").Value = "backplane-breakout-big:SAMTEC_SSM-130-L-DV"
$ws.Range("F8").Value = "SSM-130-L-DV"

$ws.Range("A9").Formula = "=A8+1"
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = "J3, J6, J7"
$ws.Range("D9").Value = "HLE-104-02-G-DV-P-TR"
$ws.Range("E9").Value = "Connector_PinHeader_2.54mm:PinHeader_2x04_P2.54mm_Vertical_SMD"
$ws.Range("F9").Value = "HLE-104-02-G-DV-P-TR"

$ws.Range("A10").Formula = "=A9+1"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "J8"
$ws.Range("D10").Value = "HSEC8-160-01-S-DV-A-K-TR"
$ws.Range("E10").Value = "backplane-breakout-big:HSEC8-160-01-S-DV-A-K-TR"
$ws.Range("F10").Value = "HSEC8-160-01-S-DV-A-K-TR"

$ws.Range("A11").Formula = "=A10+1"
$ws.Range("B11").Value = 7
$ws.Range("C11").Value = "Q1, Q4, Q5, Q6, Q9, Q10, Q13"
$ws.Range("D11").Value = "BSS84W-7-F"
$ws.Range("E11").Value = "Package_TO_SOT_SMD:SOT-323_SC-70"
$ws.Range("F11").Value = "BSS84W-7-F"

$ws.Range("A12").Formula = "=A11+1"
$ws.Range("B12").Value = 16
$ws.Range("C12").Value = "Q2, Q3, Q7, Q8, Q11, Q12, Q14, Q15, Q16, Q17, Q18, Q19, Q20, Q21, Q22, Q23"
$ws.Range("D12").Value = "DMG1012UW-7"
$ws.Range("E12").Value = "Package_TO_SOT_SMD:SOT-323_SC-70"
$ws.Range("F12").Value = "DMG1012UW-7"

$ws.Range("A13").Formula = "=A12+1"
$ws.Range("B13").Value = 16
$ws.Range("C13").Value = "R1, R2, R8, R9, R14, R15, R20, R21, R24, R25, R26, R27, R32, R33, R34, R35"
$ws.Range("D13").Value = "10k"
$ws.Range("E13").Value = "Resistor_SMD:R_0603_1608Metric"
$ws.Range("F13").Value = "RMCF0603FT100K"

$ws.Range("A14").Formula = "=A13+1"
$ws.Range("B14").Value = 23
$ws.Range("C14").Value = "R3, R4, R5, R6, R7, R10, R11, R12, R13, R16, R17, R18, R19, R22, R23, R28, R29, R30, R31, R36, R37, R38, R39"
$ws.Range("D14").Value = "1k"
$ws.Range("E14").Value = "Resistor_SMD:R_0603_1608Metric"
$ws.Range("F14").Value = "ESR03EZPF1001"

$ws.Range("A15").Formula = "=A14+1"
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = "SW1,  SW2"
$ws.Range("D15").Value = "BUS_RESET"
$ws.Range("E15").Value = "backplane-breakout-big:SW_2-1437565-9"
$ws.Range("F15").Value = "2-1437565-9"

$ws.Range("A16").Formula = "=A15+1"
$ws.Range("B16").Value = 9
$ws.Range("C16").Value = "SW3, SW4, SW5, SW6, SW7, SW8, SW9, SW10, SW11"
$ws.Range("D16").Value = "USER_SW_1"
$ws.Range("E16").Value = "backplane-breakout-big:JS102011SCQN"
$ws.Range("F16").Value = "JS102011SCQN"

$ws.Range("A17").Formula = "=A16+1"
$ws.Range("B17").Value = 120
$ws.Range("C17").Value = "TP1 thru TP120"
$ws.Range("D17").Value = "Keystone 5027"
$ws.Range("E17").Value = "backplane-breakout-big:Keystone_5027"
$ws.Range("F17").Value = "Keystone 5027"

# Column widths (bestFit sizes from target)
$ws.Columns("C").ColumnWidth = 93.140625
$ws.Columns("D").ColumnWidth = 25
$ws.Columns("E").ColumnWidth = 67.28515625
$ws.Columns("F").ColumnWidth = 26.7109375

# Selection matches source workbook state
$ws.Range("C10").Select()
